$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (case changed: ToolNaam -> Toolnaam, etc.)
$ws.Range("A1").Value = "Toolnaam"
$ws.Range("B1").Value = "Toolbeschrijving"
$ws.Range("C1").Value = "Toolimg"
$ws.Range("D1").Value = "Toollink"

# Row 2 (E-health) keeps the same description/image/Null values; re-set for safety
$ws.Range("B2").Value = "We leven allemaal tegen woordig digitaal dus waarom niet bellen met je dokter. Of inzien hoe het met je gaat. Dat allemaal digitaal bij de hand."
$ws.Range("C2").Value = "https://www.ictmagazine.nl/wp-content/uploads/2017/06/Cormel.png"
$ws.Range("D2").Value = "Null"

# Row 3: Domotica
$ws.Range("B3").Value = "Hoe werken IOT en Smart mobiles nou samen met de zorg? Vindt het hier uit"
$ws.Range("C3").Value = "https://www.duurzaammbo.nl/images/foto2/domotica.jpg"
$ws.Range("D3").Value = "Null"

# Row 4: MomoBedsense
$ws.Range("B4").Value = "Bedsensoren. Tenslotte willen we niet dat onze patienten een onrustigge nacht hebben. Hoe houden we dat in de gaten?"
$ws.Range("C4").Value = "https://www.fundis.nl/wp-content/uploads/2019/05/FundiQare_Momo-Medical_Plaat-en-Box.jpg"
$ws.Range("D4").Value = "Null"

# Row 5: Persoonsalarmering
$ws.Range("B5").Value = "Persoon kwijt of iemand met dementie. Maak gebruik van de Persoon alarm. Vaak gepaard met gps."
$ws.Range("C5").Value = "https://www.curamare.nl/images/content/page341/3ded8d9da1curamare-persoonsalarmering.jpg"
$ws.Range("D5").Value = "Null"

# Add hyperlinks to the newly-populated image-link cells, matching style of C2
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.duurzaammbo.nl/images/foto2/domotica.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.fundis.nl/wp-content/uploads/2019/05/FundiQare_Momo-Medical_Plaat-en-Box.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.curamare.nl/images/content/page341/3ded8d9da1curamare-persoonsalarmering.jpg") | Out-Null

$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("C4").Style = $ws.Range("C2").Style
$ws.Range("C5").Style = $ws.Range("C2").Style

# Update selection to C14 per diff
$ws.Range("C14").Select() | Out-Null
